$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 17 (Loads_TC007): replace the "Import button" description with the
#    Shipper-Admin-specific version of the Import File test case.
$ws.Range("B17").Value = "Validate the Import File using Shipper Admin.`n1) Enter valid user id and Password and click Login button.`n2) Observe loads are displayed in Load grid.`n3) Select one or more loads and click Import button.`n4) Upload the excel sheet .`n5) Finish the load."

# 2) Row 18 (Loads_TC008): replace with the Global-Admin-specific version.
$ws.Range("B18").Value = "Validate the Import File using Global Admin.`n1) Enter valid user id and Password and click Login button.`n2) Observe loads are displayed in Load grid.`n3) Select one or more loads and click Import button.`n4) Upload the excel sheet .`n5) Finish the load."

# 3) Add a brand new row 19 (Loads_TC009) - copy formatting/styles from row 18
#    first, then overwrite the text that differs.
$ws.Range("A18:D18").Copy($ws.Range("A19:D19"))
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

$ws.Range("A19").Value = "Loads_TC009"
$ws.Range("B19").Value = "Validate the Import File using Shipper User.`n1) Enter valid user id and Password and click Login button.`n2) Observe loads are displayed in Load grid.`n3) Select one or more loads and click Import button.`n4) Upload the excel sheet .`n5) Finish the load."
$ws.Range("C19").Value = "No"
$ws.Range("D19").Value = "Import Done Successfully"

# 4) Update the view state to match: scrolled down one row, new selection on B19.
$ws.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 16
